$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "E3"  = 16.272
    "E21" = 16.541
    "E23" = 16.55
    "E25" = 17.493
    "E53" = 16.639
    "E57" = 16.556
    "E59" = 16.305
    "E69" = 17.563
    "E79" = 17.195
    "E83" = 16.535
    "E93" = 17.103
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
